$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.0
$ws.Range("C2").Value = 74.0
$ws.Range("F2").Value = 0.03857059124412755
$ws.Range("G2").Value = 0.027246838052063738
$ws.Range("H2").Value = 0.0473148015517764
$ws.Range("I2").Value = 0.06698136254656667
